$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Round the Ost (Q2) and Nord (R2) coordinate values to nearest integer
$ws.Range("Q2").Value = 332333
$ws.Range("R2").Value = 6626960

# Remove the Starttid (Z2) and Sluttid (AB2) inline string values
$ws.Range("Z2").ClearContents()
$ws.Range("AB2").ClearContents()
